# A new daily price record was inserted into the dataset at row 434,
# shifting all subsequent rows (old 434-512) down by one (new 435-513).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 434; rows 434-512 move down to 435-513,
# and the new row inherits formatting (e.g. the date number format) from
# the row above, same as Excel's normal Insert behavior.
$ws.Rows.Item(434).Insert()

# Populate the newly inserted row 434 with the new record's data.
$ws.Range("A434").Value = 4
$ws.Range("B434").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C434").Value = "Los Lagos"
$ws.Range("D434").Value = 45209
$ws.Range("E434").Value = 10
$ws.Range("F434").Value = 100112043
$ws.Range("G434").Value = "Pepino ensalada"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 400
$ws.Range("K434").Value = 18000
$ws.Range("L434").Value = 20000
$ws.Range("M434").Value = 19000
$ws.Range("N434").Value = "$/caja 60 unidades"
$ws.Range("O434").Value = "Región de Arica y Parinacota"
$ws.Range("P434").Value = 317
$ws.Range("Q434").Value = 60
$ws.Range("R434").Value = "Hortaliza"
